$wb = $excel.ActiveWorkbook

# Update the "Port-info" header to "Port-comment" on every sheet (E2 cell)
foreach ($ws in $wb.Worksheets) {
    $ws.Range("E2").Value = "Port-comment"
}

$wsUart   = $wb.Worksheets.Item("uart")
$wsUartRx = $wb.Worksheets.Item("uart_rx")
$wsUartTx = $wb.Worksheets.Item("uart_tx")

# Sheet "uart": remove stray debug comments
$wsUart.Range("E6").ClearContents()
$wsUart.Range("E9").ClearContents()

# Sheet "uart_rx": remove stray debug comments
$wsUartRx.Range("E4").ClearContents()
$wsUartRx.Range("E7").ClearContents()
$wsUartRx.Range("E11").ClearContents()

# Sheet "uart_tx": remove/replace stray debug comments with real ones
$wsUartTx.Range("E4").ClearContents()
$wsUartTx.Range("E5").Value = "dsaf"
$wsUartTx.Range("E7").Value = "asf"
$wsUartTx.Range("E9").ClearContents()
$wsUartTx.Range("E10").Value = "fdasf"
